# Updated Logger class and Program Delete Scenario
#
# - Add a new "create" worksheet (cloned formatting from "postdata"),
#   with test data for a "create program" scenario.
# - Keep "postdata" and "putdata" worksheets (with "putdata" data updated).
# - Delete the "putProgramname" worksheet entirely.
#
# NOTE: worksheet references captured before a Copy()/rename/delete can
# become stale (they resolve by position, which shifts), so we re-fetch
# each worksheet by name right before using it.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "create" sheet by cloning "postdata" so it keeps the
#        same column widths / cell styles (date-format style, etc.) ---
$postdata = $wb.Worksheets.Item("postdata")
$null = $postdata.Copy($postdata)

$createSheet = $wb.Worksheets.Item("postdata (2)")
$createSheet.Name = "create"

# Update row 3 values for the create scenario
$createSheet.Range("A3").Value = "NinjaSparrree"
$createSheet.Range("B3").Value = "Learn Automation"
$createSheet.Range("C3").Value = "Active"

# Add row 4 with a second create-scenario entry, reusing formatting from row 3
$null = $createSheet.Range("A3:D3").Copy()
$null = $createSheet.Range("A4:D4").PasteSpecial(-4122)
$createSheet.Range("A4").Value = "NinjaDiiiee"
$createSheet.Range("B4").Value = "Learn API Test"
$createSheet.Range("C4").Value = "Active"

# --- 2. "postdata" sheet data/text itself is unchanged (NinjaSparkSDET1 / Learn API) ---
# (left as-is)

# --- 3. Update "putdata" sheet values ---
$putdata = $wb.Worksheets.Item("putdata")
$putdata.Range("A2").Value = "NinjaSparkCin221"
$putdata.Range("A3").Value = "NinjaSparkCin231"

# --- 4. Delete the "putProgramname" sheet ---
$putProgramname = $wb.Worksheets.Item("putProgramname")
$null = $putProgramname.Delete()

# --- 5. Adjust selections / active sheet to match final view state ---
$putdata = $wb.Worksheets.Item("putdata")
$null = $putdata.Activate()
$null = $putdata.Range("B7").Select()

$postdata = $wb.Worksheets.Item("postdata")
$null = $postdata.Activate()
$null = $postdata.Range("B3").Select()

$createSheet = $wb.Worksheets.Item("create")
$null = $createSheet.Activate()
$null = $createSheet.Range("B4").Select()
